$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- T8: re-style (yellow fill + full border), keep its value ---
[void]$ws.Range("U3").Copy()
$ws.Range("T8").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 11 ---
[void]$ws.Range("C3").Copy()
$ws.Range("C11").PasteSpecial(-4104)  # xlPasteAll ("0.9")

$ws.Range("H11").Formula = '="0.89"'
[void]$ws.Range("H11").Copy()
$ws.Range("H11").PasteSpecial(-4163)  # xlPasteValues -> plain text "0.89", no style change
[void]$ws.Range("C3").Copy()
$ws.Range("H11").PasteSpecial(-4122)  # xlPasteFormats -> restore style s=3 without new style entries

$ws.Range("I11").Value = 4

[void]$ws.Range("J3").Copy()
$ws.Range("J11").PasteSpecial(-4104)  # "0.00002"

[void]$ws.Range("L3").Copy()
$ws.Range("L11").PasteSpecial(-4104)  # "0.0001"

[void]$ws.Range("P3").Copy()
$ws.Range("P11").PasteSpecial(-4104)  # "0.0"

# --- Row 12 ---
[void]$ws.Range("C3").Copy()
$ws.Range("C12").PasteSpecial(-4104)  # "0.9"

$ws.Range("H12").Formula = '="0.91"'
[void]$ws.Range("H12").Copy()
$ws.Range("H12").PasteSpecial(-4163)  # xlPasteValues -> plain text "0.91", no style change
[void]$ws.Range("C3").Copy()
$ws.Range("H12").PasteSpecial(-4122)  # xlPasteFormats -> restore style s=3 without new style entries

$ws.Range("I12").Value = 4

[void]$ws.Range("J3").Copy()
$ws.Range("J12").PasteSpecial(-4104)  # "0.00002"

[void]$ws.Range("L3").Copy()
$ws.Range("L12").PasteSpecial(-4104)  # "0.0001"

[void]$ws.Range("P3").Copy()
$ws.Range("P12").PasteSpecial(-4104)  # "0.0"

# --- Row 13 ---
[void]$ws.Range("C3").Copy()
$ws.Range("C13").PasteSpecial(-4104)  # "0.9"

[void]$ws.Range("H10").Copy()
$ws.Range("H13").PasteSpecial(-4104)  # "0.925"

$ws.Range("I13").Value = 4

[void]$ws.Range("J3").Copy()
$ws.Range("J13").PasteSpecial(-4104)  # "0.00002"

# --- Row 14 ---
[void]$ws.Range("C3").Copy()
$ws.Range("C14").PasteSpecial(-4104)  # "0.9"

[void]$ws.Range("H9").Copy()
$ws.Range("H14").PasteSpecial(-4104)  # "0.875"

$ws.Range("I14").Value = 4

[void]$ws.Range("J3").Copy()
$ws.Range("J14").PasteSpecial(-4104)  # "0.00002"

# --- Selection moves to I16 ---
[void]$ws.Range("I16").Select()
